$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update numeric statistics in column B
$ws.Range("B2").Value2 = 367.0
$ws.Range("B6").Value2 = 310.0
$ws.Range("B9").Value2 = 60.0

# Fix the fragmented Polish words (missing diacritics merged back in) in column C
$words = @("Zaczął", "Pan", "nam", "podsuwać", "różne", "lektury", "które", "mają", "nas", "rozwijać", "w", "sposoby", "Ciekawi", "mnie", "jaką", "opinie", "ma", "pan", "na", "temat", "relacji", "damsko", "męskich", "czy", "jesteśmy", "odpowiednim", "wieku", "do", "angażowania", "się", "związki", "Nie", "powinniśmy", "skupić", "naszym", "rozwoju", "Przecież", "inwestowanie", "czasu", "kogoś", "z", "kim", "zapewne", "kiedyś", "rozstaniemy", "wydaje", "zaprzeczeniem", "tego", "czego", "profesor", "uczy")
for ($i = 0; $i -lt $words.Length; $i++) {
    $ws.Cells.Item(3 + $i, 3).Value2 = $words[$i]
}

# Remove the now-obsolete trailing rows (their words were merged above)
$ws.Rows("54:59").Delete()

